$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price (D) and volume-change (E) figures, and the
# re-ordered Cosmos/Toncoin and ARBITRUM/Aave/FTXToken rows, as of
# the Thu Dec 7 22:59:40 UTC 2023 GitHub Actions refresh.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.262.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.357.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.63%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.649"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.30"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.62%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.457"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0954"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.710.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.56%  "
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.839"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.362.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.236.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0980"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +18.06%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.61%  "
$ws.Range("E32").Value = "  -6.74%  "
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.82%  "
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.20%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "18.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.56%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.73%  "
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0948"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.444.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.583.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.91%  "
